$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("source_kb_shots.csv")
$ws.Activate()

# Append the new "Variables supplémentaires" rows (43-52) below the
# existing list in columns B (name) / C (description).
$newRows = @(
    @("boo_premier_shot_qt",    "Booléen premier shot de KB du QT"),
    @("boo_dernier_shot_qt",    "Booléen dernier shot de KB du QT"),
    @("boo_premier_shot_match", "Booléen premier shot de KB du match"),
    @("boo_dernier_shot_match", "Booléen dernier shot de KB du match"),
    @("temps_dernier_shot",     "Temps en seconde depuis le dernier shot (si premier shot du QT on met 12*60=720 secondes"),
    @("temps_prochain_shot",    "Temps en seconde avat le prochain shot (si dernier shot du QT on met 12*60=720 secondes"),
    @("nb_shot_qt",             "Nombre de shot pris par KB au total dans le QT"),
    @("intensite_shot_qt",      "Intensité des shots pris par KB sur ce QT (nombre de shot dans le QT dvisié par la durée du QT)"),
    @("nb_shot_match",          "Nombre de shot pris par KB au total dans le match"),
    @("intensite_shot_match",   "Intensité des shots pris par KB sur ce match (nombre de shot dans le match dvisié par la durée du match)")
)

$startRow = 43
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 2).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 3).Value = $newRows[$i][1]
}

# Reflect the selection recorded in the saved file (the view also scrolls
# so row 34 is at the top, but that scroll position isn't a separately
# addressable/persisted property in this host).
$excel.ActiveWindow.ScrollRow = 34
$ws.Range("C54").Select()
